# Apply "Added Mid-Game attacks" updates to the project planning sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Implement Core Mechanics (Fight, Run) - Current Est increased 15 -> 20
$ws.Range("C3").Value = 20

# > Implement Items - Current Est increased 10 -> 15
$ws.Range("C4").Value = 15

# > Implement Shield Mechanic - Current Est increased 8 -> 12
$ws.Range("C6").Value = 12

# Writing the Story - Effort increased 1 -> 2
$ws.Range("D13").Value = 2

# Implement attacks (Mid-Game attacks) - Effort increased 5 -> 8
$ws.Range("D14").Value = 8

# Move the active selection to E14
$ws.Range("E14").Select()
